# Update the RouteWork row (row 2) with the 20th June 2022 values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RouteWorkStartDate (P2) / AlternateGenerationDate (CA2): 10-May-2022 -> 17-Jun-2022
$ws.Range("P2").Value = 44729
$ws.Range("CA2").Value = 44729

# RouteWorkReadyTime (S2): 00:00 -> 07:00
$ws.Range("S2").Value = 0.29166666666666669

# RouteWorkScheduledEndTime (T2): 23:15 -> 06:45
$ws.Range("T2").Value = 0.28125

# FirstGenerationTime (CB2): 00:05 -> 07:00
$ws.Range("CB2").Value = 0.29166666666666669

# Move the sheet's current selection to CA2, matching the saved workbook state.
$ws.Range("CA2").Select()
